# Update "想去人数" (F column) figures across sheets, matching output
# regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 9477
$ws1.Range("F5").Value = 614
$ws1.Range("F6").Value = 161
$ws1.Range("F7").Value = 288
$ws1.Range("F8").Value = 354
$ws1.Range("F9").Value = 410
$ws1.Range("F11").Value = 187
$ws1.Range("F12").Value = 18
$ws1.Range("F14").Value = 12122
$ws1.Range("F25").Value = 2099
$ws1.Range("F30").Value = 1009
$ws1.Range("F32").Value = 3649
$ws1.Range("F33").Value = 564
$ws1.Range("F35").Value = 3060
$ws1.Range("F37").Value = 1324
$ws1.Range("F40").Value = 19
$ws1.Range("F43").Value = 535
$ws1.Range("F46").Value = 227
$ws1.Range("F47").Value = 112
$ws1.Range("F49").Value = 142

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 16
$ws2.Range("F19").Value = 10

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 9477
$ws4.Range("F7").Value = 614
$ws4.Range("F9").Value = 161
$ws4.Range("F10").Value = 288
$ws4.Range("F11").Value = 354
$ws4.Range("F12").Value = 410
$ws4.Range("F14").Value = 187
$ws4.Range("F16").Value = 12122
$ws4.Range("F26").Value = 2099
$ws4.Range("F31").Value = 1009
$ws4.Range("F33").Value = 3649
$ws4.Range("F34").Value = 564
$ws4.Range("F36").Value = 3060
$ws4.Range("F37").Value = 1324
$ws4.Range("F43").Value = 535
$ws4.Range("F46").Value = 227
$ws4.Range("F47").Value = 112
$ws4.Range("F49").Value = 142
